$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Payton Pritchard", "PG,SG", "Boston Celtics"),
    @("Jalen Green", "PG,SG", "Houston Rockets"),
    @("Ty Jerome", "PG,SG", "Cleveland Cavaliers"),
    @("Brandin Podziemski", "PG,SG", "Golden State Warriors"),
    @("Deni Avdija", "SF,PF", "Portland Trail Blazers"),
    @("Paolo Banchero", "SF,PF", "Orlando Magic"),
    @("Jaylen Brown", "SG,SF", "Boston Celtics"),
    @("Kyle Filipowski", "PF,C", "Utah Jazz"),
    @("Jakob Poeltl", "C", "Utah Jazz"),
    @("Pascal Siakam", "SF,PF,C", "Indiana Pacers"),
    @("Nikola Jokic", "C", "Denver Nuggets"),
    @("Andre Drummond", "C", "Philadelphia 76ers"),
    @("Chris Paul", "PG", "San Antonio Spurs"),
    @("Russell Westbrook", "PG,SG", "Denver Nuggets"),
    @("Stephon Castle", "PG,SG", "San Antonio Spurs"),
    @("Rudy Gobert", "C", "Minnesota Timberwolves"),
    @("Jerami Grant", "SF,PF", "Portland Trail Blazers"),
    @("Chet Holmgren", "PF,C", "Oklahoma City Thunder")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
